# QA Round 2: deep quality optimization - compliance, diversification, UX improvements
#
# 1) FabyJourney (sheet1): a few line-text tweaks.
# 2) cumcontrol -> split into "cumcontrol1" (edited in place, renamed) and a brand
#    new "cumcontrol2" sheet (inserted right before "dickpic") carrying a second
#    diversified set of delay/sync/edge variants. The new sheet is produced by
#    duplicating "cumcontrol" (so it inherits the same layout/column widths/
#    cell styles) and then overwriting its text.
# 3) "dickpic" and "boosters" keep their original content; they just shift right
#    in the tab order to make room for the new "cumcontrol2" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) FabyJourney tweaks
# ---------------------------------------------------------------------------
$faby = $wb.Worksheets.Item("FabyJourney")
$faby.Range("B4").Value  = "finish with me gostoso"
$faby.Range("B8").Value  = "oh fuck"
$faby.Range("B10").Value = "gimme a minute"

# ---------------------------------------------------------------------------
# 2) Duplicate "cumcontrol" right after itself -> becomes "cumcontrol2"
#    (placed before "dickpic", inherits formatting/column widths/styles)
# ---------------------------------------------------------------------------
$cc1 = $wb.Worksheets.Item("cumcontrol")
$cc1.Copy($null, $cc1)
$cc2 = $wb.Worksheets.Item("cumcontrol (2)")

# ---------------------------------------------------------------------------
# cumcontrol -> cumcontrol1 (rename + in-place text edits)
# ---------------------------------------------------------------------------
$cc1.Range("B2").Value = "patience... what's coming is worth every second of waiting"

$cc1.Range("B3").Value = "hold it for me... I have years of experience and this next one is my best work"
$cc1.Range("C3").Value = "DELAY. Send PPV."

$cc1.Range("B4").Value = "I want to feel you let go while I do the same... watch this first"
$cc1.Range("C4").Value = "SYNC variant. Send PPV."

$cc1.Range("B5").Value = "now we go together amor... I've been holding back too. open this 😏"
$cc1.Range("C5").Value = "SYNC. Send PPV."

$cc1.Range("B6").Value = "a man who can wait gets rewarded gostoso... trust me on that"

$cc1.Range("B7").Value = "I can tell you're close... not yet amor, I know what I'm doing"
$cc1.Range("C7").Value = "CONTROL."

$cc1.Name = "cumcontrol1"

# ---------------------------------------------------------------------------
# cumcontrol (2) -> cumcontrol2 (new diversified delay/sync/edge variants)
# ---------------------------------------------------------------------------
$cc2.Range("A2").Value = "delay2"
$cc2.Range("B2").Value = "save it for this last one amor, I promise you it's going to be worth it 😏"
$cc2.Range("C2").Value = "DELAY variant."

$cc2.Range("A3").Value = "delay1"
$cc2.Range("B3").Value = "one more for you before we're done... this is the one I'm most proud of"
$cc2.Range("C3").Value = "DELAY. Send PPV."

$cc2.Range("A4").Value = "sync2"
$cc2.Range("B4").Value = "I'm ready when you are... but see this first"
$cc2.Range("C4").Value = "SYNC variant."

$cc2.Range("A5").Value = "sync1"
$cc2.Range("B5").Value = "okay amor... let's both let go right now. open this 😏"
$cc2.Range("C5").Value = "SYNC. Send PPV."

$cc2.Range("A6").Value = "edge2"
$cc2.Range("B6").Value = "not yet... a little more anticipation makes it so much better, trust me"
$cc2.Range("C6").Value = "EDGE variant."

$cc2.Range("A7").Value = "edge1"
$cc2.Range("B7").Value = "slow down for me gostoso... I know exactly when to let you go"
$cc2.Range("C7").Value = "CONTROL."

# Note: column D ("*Guidelines") cells are already blank in the duplicated
# sheet (inherited from "cumcontrol"), so nothing further to clear there.

$cc2.Name = "cumcontrol2"
